$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.255.18'
$ws.Range('E2').Value = '  -1.88%  '
$ws.Range('D3').Value = '2.236.63'
$ws.Range('E3').Value = '  -1.81%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.92'
$ws.Range('E5').Value = '  -1.05%  '
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '63.38'
$ws.Range('E7').Value = '  -1.15%  '
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.438'
$ws.Range('E9').Value = '  +0.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0945'
$ws.Range('E10').Value = '  -8.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '56.20'
$ws.Range('E11').Value = '  -0.76%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '27.43'
$ws.Range('E12').Value = '  +5.81%  '
$ws.Range('E13').Value = '  -1.82%  '
$ws.Range('D14').Value = '2.568.49'
$ws.Range('E14').Value = '  -1.84%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.16'
$ws.Range('E15').Value = '  -4.38%  '
$ws.Range('E16').Value = '  +0.92%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.821'
$ws.Range('E17').Value = '  -1.07%  '
$ws.Range('D18').Value = '2.247.21'
$ws.Range('E18').Value = '  -1.27%  '
$ws.Range('D19').Value = '43.100.68'
$ws.Range('E19').Value = '  -2.10%  '
$ws.Range('D20').Value = '0.0₃0960'
$ws.Range('E20').Value = '  -8.77%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.75'
$ws.Range('E21').Value = '  -1.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.07'
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '245.75'
$ws.Range('E23').Value = '  -4.77%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.69'
$ws.Range('E25').Value = '  +29.76%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.42'
$ws.Range('E26').Value = '  -2.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.29'
$ws.Range('E27').Value = '  -1.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.68'
$ws.Range('E28').Value = '  -3.93%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '173.07'
$ws.Range('E29').Value = '  +0.83%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '21.50'
$ws.Range('E30').Value = '  +2.36%  '
$ws.Range('E31').Value = '  -7.48%  '
$ws.Range('E32').Value = '  -1.07%  '
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('E34').Value = '  +3.49%  '
$ws.Range('E35').Value = '  -1.65%  '
$ws.Range('E36').Value = '  -2.03%  '
$ws.Range('E37').Value = '  -7.31%  '
$ws.Range('E38').Value = '  -7.59%  '
$ws.Range('E40').Value = '  -0.50%  '
$ws.Range('E41').Value = '  +0.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.61'
$ws.Range('E42').Value = '  +1.64%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.43'
$ws.Range('E43').Value = '  -1.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.88'
$ws.Range('E44').Value = '  -4.51%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0941'
$ws.Range('E45').Value = '  -3.18%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '96.05'
$ws.Range('E46').Value = '  -2.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.17'
$ws.Range('E47').Value = '  -2.50%  '
$ws.Range('E48').Value = '  -1.66%  '
$ws.Range('D49').Value = '1.435.53'
$ws.Range('E49').Value = '  -2.36%  '
$ws.Range('E50').Value = '  +1.95%  '
$ws.Range('E51').Value = '  +0.41%  '
